$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "sotl"/"file_name"/"impact"/"Impact"/... rows and the stray D column,
# replacing them with the new glossary content (A:C, 6 rows).
$ws.Range("A1:D3").ClearContents()

# Header row and the existing "SoTL" row stay row-major (already present before this edit).
$ws.Range("A1").Value = "term"
$ws.Range("B1").Value = "short_def"
$ws.Range("C1").Value = "long_def"

$ws.Range("A2").Value = "SoTL"
$ws.Range("B2").Value = "A scholarly approach to teaching."
$ws.Range("C2").Value = "SoTL, or the Scholarship of Teaching and Learning, is a scholarly approach to teaching that involves systematically examining one's own teaching practices and student learning to improve both, often with the goal of making findings public to contribute to the wider teaching community"

# New terms: type all of column A first ...
$ws.Range("A3").Value = "Pedagogy"
$ws.Range("A4").Value = "Active learning"
$ws.Range("A5").Value = "Engagement"
$ws.Range("A6").Value = "Autonomy"

# ... then go back and fill in the long definition (C) before the short definition (B)
# for each new row.
$ws.Range("C3").Value = "The method and practice of teaching, especially as an academic subject or theoretical concept."
$ws.Range("B3").Value = "The method and practice of teaching."

$ws.Range("C4").Value = "Active learning is a method of learning in which students are actively or experientially involved in the learning process."
$ws.Range("B4").Value = "Students actively involved in the learning process."

$ws.Range("C5").Value = "Psychological investiment in learning."
$ws.Range("B5").Value = "Psychological investiment in learning."

$ws.Range("C6").Value = "`nAutonomy means the ability to take control of one's own learning, independently or in collaboration with others."
$ws.Range("B6").Value = "Taking control of one's own learning."

$ws.Range("C6").WrapText = $true
$ws.Rows.Item(6).EntireRow.AutoFit()

$ws.Range("B7").Select()
